# Commit: Tue, May 19, 2020 11:05:38 AM
#
# 1) Slide 16 ("PLENARY- COMPLETE THE MISSING GAPS") has a 2-column table
#    (3rd shape on the slide, a graphicFrame) whose table style was changed
#    from the default generated style {624EDA45-2844-48E4-8D4C-FB826F15C6FB}
#    to the built-in PowerPoint table style {142CDE50-73C3-4405-BFC5-C3BE210F9A24}.
#
# 2) The deck's "Integral" theme (ppt/theme/theme1.xml, used by the slide
#    master/all slides) and its "Office Theme" (ppt/theme/theme2.xml, used
#    only by the notes master) had their color palettes swapped. Apply the
#    "Office Theme" palette to the presentation's live theme color scheme.

$p = $ppt.ActivePresentation

# --- 1) Table style on slide 16 ---------------------------------------
$slide = $p.Slides.Item(16)
$tableShape = $slide.Shapes.Item(3)
$table = $tableShape.Table
$table.ApplyStyle("{142CDE50-73C3-4405-BFC5-C3BE210F9A24}")

# --- 2) Theme color palette swap (Integral -> Office Theme colors) ----
$themeColors = $slide.ThemeColorScheme

# Office theme palette values, in COM RGB (0x00BBGGRR) form:
# dk1, lt1, dk2, lt2, accent1..accent6, hlink, folHlink
$officeColors = @(
    0x000000,  # dk1      000000
    0xFFFFFF,  # lt1      FFFFFF
    0x6A5444,  # dk2      44546A
    0xE6E6E7,  # lt2      E7E6E6
    0xD59B5B,  # accent1  5B9BD5
    0x317DED,  # accent2  ED7D31
    0xA5A5A5,  # accent3  A5A5A5
    0x00C0FF,  # accent4  FFC000
    0xC47244,  # accent5  4472C4
    0x47AD70,  # accent6  70AD47
    0xC16305,  # hlink    0563C1
    0x724F95   # folHlink 954F72
)

for ($i = 1; $i -le $officeColors.Count; $i++) {
    $themeColors.Colors($i).RGB = $officeColors[$i - 1]
}
